$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Items")

# Add a new row (row 6) for the "unpublish" item, cloning the formatting/
# structure of the last existing data row (row 5) and then overwriting the
# two cells that differ: the item ID and the action.
$ws.Range("A5:R5").Copy($ws.Range("A6:R6"))
$ws.Range("A6").Value = "ITM-1213-3316-0005"
$ws.Range("C6").Value = "unpublish"

# Extend the Action column's data validation list so the whole column
# (minus the header) allows the new "unpublish" choice, matching the
# widened sqref down to the sheet's row limit.
$ws.Range("C2:C1048576").Validation.Delete()
$ws.Range("C2:C1048576").Validation.Add(3, 1, 1, '"-,update,review,publish,unpublish"')
$ws.Range("C2:C1048576").Validation.IgnoreBlank = $false

# Move the active selection to C9, as recorded in the saved view state.
$ws.Activate()
$ws.Range("C9").Select()
